# WCAG_1.4.3_MinimumContrast test case workbook update
# - Add a "Level Classification: AA" field to the summary block (E12/F12)
# - Give the big Step 1 / Step 2 rows (14-15) a taller, vertically centered layout
# - Re-point the active selection at E12 and set the page to Portrait orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Level Classification" / "AA" fields (row 12) -------------------
$ws.Range("E12").Value = "Level Classification"
$ws.Range("F12").Value = "AA"

# Match the look of the other field labels in column E (bold, shaded, bordered)
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Rows 14 & 15: taller rows, vertically centered content ---------------
$ws.Range("A14:K14").VerticalAlignment = -4108
$ws.Range("A15:K15").VerticalAlignment = -4108

$ws.Rows.Item(14).RowHeight = 247.5
$ws.Rows.Item(15).RowHeight = 180

# --- View state: select E12, switch the sheet to portrait page orientation -
$ws.Range("E12").Select() | Out-Null
$ws.PageSetup.Orientation = 1

Write-Output "applied updates"
